$wb = $excel.ActiveWorkbook

# --- select A3 on currency_conversions (non-active-sheet selection) ---
$convWs = $wb.Worksheets.Item("currency_conversions")
[void]$convWs.Range("A3").Select()

# --- add new sheet "currency_movements" after the last existing sheet ---
$lastWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastWs)
$newWs.Name = "currency_movements"

# --- bring over the existing formatting (bold header, date col, plain col)
#     from currency_conversions so we reuse the workbook's existing styles ---
$convWs.Range("B1").Copy()
$newWs.Range("A1:E1").PasteSpecial(-4122)   # xlPasteFormats -> bold header style

$convWs.Range("A2").Copy()
$newWs.Range("A2:B2").PasteSpecial(-4122)   # xlPasteFormats -> date style

$convWs.Range("B2:D2").Copy()
$newWs.Range("C2:E2").PasteSpecial(-4122)   # xlPasteFormats -> plain style

$convWs.Range("A2").Copy()
$newWs.Range("A3:A5").PasteSpecial(-4122)   # xlPasteFormats -> date style for remaining dates
$newWs.Range("B3").PasteSpecial(-4122)      # xlPasteFormats -> leave B3 blank but date-styled

# --- header row ---
$newWs.Range("A1").Value = "date"
$newWs.Range("B1").Value = "buy_date"
$newWs.Range("C1").Value = "amount"
$newWs.Range("D1").Value = "currency"
$newWs.Range("E1").Value = "comment"

# --- data rows ---
$d1 = Get-Date -Year 2022 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$d2 = Get-Date -Year 2022 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0
$d3 = Get-Date -Year 2022 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0

$newWs.Range("A2").Value = $d1
$newWs.Range("B2").Value = $d1
$newWs.Range("C2").Value = 1
$newWs.Range("D2").Value = "USD"

$newWs.Range("A3").Value = $d2
$newWs.Range("C3").Value = -1
$newWs.Range("D3").Value = "USD"

$newWs.Range("A4").Value = $d2
$newWs.Range("C4").Value = 1
$newWs.Range("D4").Value = "EUR"

$newWs.Range("A5").Value = $d3
$newWs.Range("C5").Value = -1
$newWs.Range("D5").Value = "EUR"

# --- activate the new sheet last so it becomes the selected/active tab ---
[void]$newWs.Activate()
[void]$newWs.Range("A1").Select()
